# Atualização de bases das ligas, do dia: 30-03-2024 às 19:32
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 91 and 92: the two fixtures swap their full data (everything except
# the running index in column A, which stays 89 / 90 respectively).
# ---------------------------------------------------------------------------
$ws.Range("B91").Value = 6782568
$ws.Range("F91").Value = "Sporting San Jose"
$ws.Range("G91").Value = "AD Guanacasteca"
$ws.Range("H91").Value = 1
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = "D"
$ws.Range("K91").Value = 1.909
$ws.Range("L91").Value = 3.6
$ws.Range("M91").Value = 3.3
$ws.Range("N91").Value = 2
$ws.Range("O91").Value = 3.6
$ws.Range("P91").Value = 3.1
$ws.Range("Q91").Value = -0.5
$ws.Range("R91").Value = 2
$ws.Range("S91").Value = 1.8
$ws.Range("T91").Value = 2.5
$ws.Range("U91").Value = 1.825
$ws.Range("V91").Value = 1.975
$ws.Range("W91").Value = -1
$ws.Range("X91").Value = 2.6
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.8
$ws.Range("AB91").Value = -1
$ws.Range("AC91").Value = 0.9750000000000001

$ws.Range("B92").Value = 6782566
$ws.Range("F92").Value = "Cartagines"
$ws.Range("G92").Value = "Deportivo Saprissa"
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 4
$ws.Range("J92").Value = "A"
$ws.Range("K92").Value = 3.2
$ws.Range("L92").Value = 3.4
$ws.Range("M92").Value = 2
$ws.Range("N92").Value = 2.9
$ws.Range("O92").Value = 3.5
$ws.Range("P92").Value = 2.15
$ws.Range("Q92").Value = 0.25
$ws.Range("R92").Value = 1.875
$ws.Range("S92").Value = 1.925
$ws.Range("T92").Value = 3
$ws.Range("U92").Value = 1.975
$ws.Range("V92").Value = 1.825
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = 1.15
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = 0.925
$ws.Range("AB92").Value = 0.9750000000000001
$ws.Range("AC92").Value = -1

# ---------------------------------------------------------------------------
# Live-odds refresh for the still-upcoming fixtures (rows 217-221). The
# fixture that used to sit in row 217 (id 7623941) is gone from the feed, the
# remaining ones had their odds / kick-off time refreshed, and the row that
# is now a stale duplicate (old row 222) is removed at the end.
# ---------------------------------------------------------------------------
$ws.Range("B217").Value = 8022822
$ws.Range("E217").Value = 45382.75
$ws.Range("F217").Value = "Puntarenas"
$ws.Range("G217").Value = "Municipal Liberia"
$ws.Range("K217").Value = 2.2
$ws.Range("L217").Value = 3.3
$ws.Range("M217").Value = 2.875
$ws.Range("N217").Value = 2.625
$ws.Range("O217").Value = 3.25
$ws.Range("P217").Value = 2.4
$ws.Range("Q217").Value = 0
$ws.Range("R217").Value = 2
$ws.Range("S217").Value = 1.8
$ws.Range("T217").Value = 2.5
$ws.Range("U217").Value = 1.85
$ws.Range("V217").Value = 1.95

$ws.Range("N218").Value = 1.75
$ws.Range("O218").Value = 3.6
$ws.Range("P218").Value = 4
$ws.Range("R218").Value = 1.775
$ws.Range("S218").Value = 2.025
$ws.Range("U218").Value = 1.95
$ws.Range("V218").Value = 1.85

$ws.Range("B219").Value = 7623943
$ws.Range("E219").Value = 45382.85416666666
$ws.Range("F219").Value = "Sporting San Jose"
$ws.Range("G219").Value = "AD Guanacasteca"
$ws.Range("K219").Value = 2.1
$ws.Range("L219").Value = 3.2
$ws.Range("M219").Value = 3.2
$ws.Range("N219").Value = 2
$ws.Range("O219").Value = 3.2
$ws.Range("P219").Value = 3.4
$ws.Range("Q219").Value = -0.25
$ws.Range("R219").Value = 1.775
$ws.Range("S219").Value = 2.025
$ws.Range("T219").Value = 2.25
$ws.Range("U219").Value = 1.85
$ws.Range("V219").Value = 1.95

$ws.Range("B220").Value = 7623939
$ws.Range("F220").Value = "Alajuelense"
$ws.Range("G220").Value = "AD San Carlos"
$ws.Range("R220").Value = 2
$ws.Range("S220").Value = 1.8
$ws.Range("T220").Value = 2.75
$ws.Range("U220").Value = 1.95
$ws.Range("V220").Value = 1.85

$ws.Range("B221").Value = 7623940
$ws.Range("E221").Value = 45383.75
$ws.Range("F221").Value = "AD Grecia"
$ws.Range("G221").Value = "Cartagines"
$ws.Range("K221").Value = 3
$ws.Range("L221").Value = 3.4
$ws.Range("M221").Value = 2.1
$ws.Range("N221").Value = 3.2
$ws.Range("O221").Value = 3.4
$ws.Range("P221").Value = 2
$ws.Range("Q221").Value = 0.25
$ws.Range("R221").Value = 2.025
$ws.Range("S221").Value = 1.775
$ws.Range("T221").Value = 2.5
$ws.Range("U221").Value = 1.85
$ws.Range("V221").Value = 1.95

# Drop the now-duplicated trailing row (old row 222 / match id 7623940 moved
# into row 221 above) so the sheet has 220 data rows again.
$ws.Rows.Item(222).Delete()
